$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  Name="RUIZ CHIROQUE CLAUDIA JUDITH";         Total=193},
    @{Row=3;  Name="FABIANA REBECA ARRUNATEGUI SILUPU";    Total=171},
    @{Row=4;  Name="LLENQUE ANTON HELEN JOHANA";            Total=169},
    @{Row=5;  Name="BANCAYAN FIESTA DILVER HUMBERTO";       Total=169},
    @{Row=6;  Name="TEMOCHE ECHE URSULA YESSENIA";          Total=165},
    @{Row=7;  Name="GONZALES FIESTAS MARIA MARIBEL";        Total=161},
    @{Row=8;  Name="ANTON INGA FATIMA DEL ROSARIO";         Total=160},
    @{Row=9;  Name="BAUTISTA CHAVESTA ERICKA MEDALIT";      Total=156},
    @{Row=10; Name="PINTADO CHASQUERO ESTEFANY";            Total=156},
    @{Row=11; Name="VELASCO PEÑA KAREN ARELLYS";            Total=149},
    @{Row=12; Name="HERNANDEZ CARNERO ARTURO SEBASTIAN";    Total=141},
    @{Row=13; Name="MONDRAGON NONAJULCA MARISOL";           Total=137},
    @{Row=14; Name="FLORES SILUPU MARY CARMEN";             Total=130},
    @{Row=15; Name="ORDINOLA JIBAJA JOSE ALBERTO";          Total=122},
    @{Row=16; Name="MORENO YANAYACO NAYLA GUADALUPE";       Total=100},
    @{Row=17; Name="CASTRO ESTRADA CINTHIA PATRICIA";       Total=95}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.Name
    $ws.Cells.Item($item.Row, 2).Value = $item.Total
}
